$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-28 05:48:11'
$ws.Range('N2').Value = '0.0 °C 5:12 TU'
$ws.Range('E3').Value = '2026-02-28 05:48:14'
$ws.Range('H3').Value = "'83%"
$ws.Range('N3').Value = '-2.8 °C 5:24 TU'
$ws.Range('O3').Value = '-0.4 °C'
$ws.Range('E4').Value = '2026-02-28 05:48:16'
$ws.Range('H4').Value = "'95%"
$ws.Range('M4').Value = '8.7 °C 5:14 TU'
$ws.Range('O4').Value = '7.9 °C'
$ws.Range('E5').Value = '2026-02-28 05:48:18'
$ws.Range('N5').Value = '-0.9 °C 5:28 TU'
$ws.Range('O5').Value = '-0.2 °C'
$ws.Range('E6').Value = '2026-02-28 05:48:21'
$ws.Range('E7').Value = '2026-02-28 05:48:23'
$ws.Range('L7').Value = '19.8 km/h - 78º 5:17 TU'
$ws.Range('M7').Value = '12.6 °C 5:26 TU'
$ws.Range('E8').Value = '2026-02-28 05:48:25'
$ws.Range('H8').Value = "'97%"
$ws.Range('L8').Value = '21.6 km/h - 67º 5:20 TU'
$ws.Range('M8').Value = '8.8 °C 5:29 TU'
$ws.Range('O8').Value = '8.5 °C'
$ws.Range('E9').Value = '2026-02-28 05:48:28'
$ws.Range('O9').Value = '7.4 °C'
$ws.Range('E10').Value = '2026-02-28 05:48:30'
$ws.Range('M10').Value = '9.0 °C 5:25 TU'
$ws.Range('O10').Value = '7.8 °C'
$ws.Range('E11').Value = '2026-02-28 05:48:32'
$ws.Range('H11').Value = "'90%"
$ws.Range('N11').Value = '1.4 °C 5:29 TU'
$ws.Range('O11').Value = '3.4 °C'
$ws.Range('E12').Value = '2026-02-28 05:48:35'
$ws.Range('E13').Value = '2026-02-28 05:48:37'
$ws.Range('J13').Value = '1026.3 hPa'
$ws.Range('K13').Value = '-0.1 MJ/m2'
$ws.Range('N13').Value = '-1.0 °C 5:13 TU'
$ws.Range('O13').Value = '1.1 °C'
$ws.Range('E14').Value = '2026-02-28 05:48:39'
$ws.Range('H14').Value = "'98%"
$ws.Range('L14').Value = '19.8 km/h - 59º 5:27 TU'
$ws.Range('M14').Value = '12.3 °C 5:29 TU'
$ws.Range('O14').Value = '10.4 °C'
$ws.Range('E15').Value = '2026-02-28 05:48:42'
$ws.Range('E16').Value = '2026-02-28 05:48:43'
$ws.Range('H16').Value = "'60%"
$ws.Range('K16').Value = '-0.1 MJ/m2'
$ws.Range('O16').Value = '-0.7 °C'
$ws.Range('E17').Value = '2026-02-28 05:48:46'
$ws.Range('N17').Value = '3.5 °C 5:18 TU'
$ws.Range('O17').Value = '4.6 °C'
$ws.Range('E18').Value = '2026-02-28 05:48:48'
$ws.Range('H18').Value = "'99%"
$ws.Range('O18').Value = '8.3 °C'
$ws.Range('E19').Value = '2026-02-28 05:48:51'
$ws.Range('H19').Value = "'66%"
$ws.Range('L19').Value = '29.2 km/h - 57º 5:25 TU'
$ws.Range('E20').Value = '2026-02-28 05:48:53'
$ws.Range('H20').Value = "'37%"
$ws.Range('O20').Value = '0.0 °C'
$ws.Range('E21').Value = '2026-02-28 05:48:55'
$ws.Range('J21').Value = '1024.1 hPa'
$ws.Range('N21').Value = '3.2 °C 5:05 TU'
$ws.Range('O21').Value = '5.1 °C'
$ws.Range('E22').Value = '2026-02-28 05:48:58'
$ws.Range('H22').Value = "'58%"
$ws.Range('L22').Value = '26.6 km/h - 118º 5:07 TU'
$ws.Range('E23').Value = '2026-02-28 05:49:00'
$ws.Range('H23').Value = "'68%"
$ws.Range('E24').Value = '2026-02-28 05:49:02'
$ws.Range('O24').Value = '6.1 °C'
$ws.Range('E25').Value = '2026-02-28 05:49:04'
$ws.Range('N25').Value = '-0.8 °C 5:07 TU'
$ws.Range('O25').Value = '0.8 °C'
$ws.Range('E26').Value = '2026-02-28 05:49:07'
$ws.Range('H26').Value = "'69%"
$ws.Range('N26').Value = '3.7 °C 5:29 TU'
$ws.Range('E27').Value = '2026-02-28 05:49:09'
$ws.Range('H27').Value = "'35%"
$ws.Range('N27').Value = '0.4 °C 5:16 TU'
$ws.Range('O27').Value = '2.2 °C'
$ws.Range('E28').Value = '2026-02-28 05:49:12'
$ws.Range('E29').Value = '2026-02-28 05:49:14'
$ws.Range('E30').Value = '2026-02-28 05:49:16'
$ws.Range('E31').Value = '2026-02-28 05:49:19'
$ws.Range('H31').Value = "'93%"
$ws.Range('N31').Value = '9.5 °C 5:12 TU'
$ws.Range('E32').Value = '2026-02-28 05:49:21'
$ws.Range('E33').Value = '2026-02-28 05:49:23'
$ws.Range('H33').Value = "'72%"
$ws.Range('O33').Value = '5.0 °C'
$ws.Range('E34').Value = '2026-02-28 05:49:26'
$ws.Range('H34').Value = "'67%"
$ws.Range('L34').Value = '11.5 km/h - 17º 5:25 TU'
$ws.Range('N34').Value = '-1.5 °C 5:19 TU'
$ws.Range('O34').Value = '-0.1 °C'
$ws.Range('E35').Value = '2026-02-28 05:49:28'
$ws.Range('H35').Value = "'83%"
$ws.Range('N35').Value = '5.4 °C 5:29 TU'
$ws.Range('O35').Value = '6.7 °C'
$ws.Range('E36').Value = '2026-02-28 05:49:30'
$ws.Range('L36').Value = '20.5 km/h - 344º 5:18 TU'
$ws.Range('N36').Value = '8.8 °C 5:29 TU'
$ws.Range('E37').Value = '2026-02-28 05:49:33'
$ws.Range('N37').Value = '3.8 °C 5:29 TU'
$ws.Range('O37').Value = '4.6 °C'
$ws.Range('E38').Value = '2026-02-28 05:49:35'
$ws.Range('M38').Value = '9.9 °C 5:27 TU'
$ws.Range('O38').Value = '9.1 °C'
$ws.Range('E39').Value = '2026-02-28 05:49:37'
$ws.Range('H39').Value = "'48%"
$ws.Range('O39').Value = '0.1 °C'
$ws.Range('E40').Value = '2026-02-28 05:49:39'
$ws.Range('J40').Value = '1024.7 hPa'
$ws.Range('N40').Value = '1.8 °C 5:29 TU'
$ws.Range('O40').Value = '3.4 °C'
$ws.Range('E41').Value = '2026-02-28 05:49:42'
$ws.Range('M41').Value = '12.0 °C 5:24 TU'
$ws.Range('E42').Value = '2026-02-28 05:49:44'
$ws.Range('E43').Value = '2026-02-28 05:49:46'
$ws.Range('N43').Value = '3.0 °C 5:20 TU'
$ws.Range('O43').Value = '3.8 °C'
$ws.Range('E44').Value = '2026-02-28 05:49:49'
$ws.Range('H44').Value = "'93%"
$ws.Range('L44').Value = '20.5 km/h - 18º 5:11 TU'
$ws.Range('E45').Value = '2026-02-28 05:49:51'
$ws.Range('J45').Value = '1024.0 hPa'
$ws.Range('N45').Value = '5.1 °C 5:27 TU'
$ws.Range('O45').Value = '6.8 °C'
$ws.Range('E46').Value = '2026-02-28 05:49:54'
$ws.Range('H46').Value = "'79%"
